# Updated cryptos list on Sun Feb 26 21:36:09 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the cryptos
# table on rows 2-51 with newly scraped values.
#
# Many "Price" values are plain decimal-looking strings (e.g. "0.9985",
# "53.00"). Excel's Range.Value setter auto-coerces such strings to
# numbers, which would both change the stored type (text -> number) and
# silently mangle the literal text (trailing zeros lost, floating point
# noise introduced). To keep these as literal text -- matching the
# original workbook's inline-string cells -- we flip the cell to a text
# number format before writing the value, then restore the cell's style
# to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.551.68"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "1.641.51"
$ws.Range("E3").Value = "  +4.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3681"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.286"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08194"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.677"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.484"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").Value = "1.640.99"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06955"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.605"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9981"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "23.559.28"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.131"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.67%  "
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.333"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.425"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.860"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.51%  "
$ws.Range("D33").Value = "1.816.06"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9780"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.54%  "
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07494"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.242"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2546"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.84%  "
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6655"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.371"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.045"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9972"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08073"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.216"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.64%  "
